# Insert a new weekly record at the top of the "Poroto verde" /
# "Feria Lagunitas de Puerto Montt" block (row 169), pushing the
# existing rows 169-179 down to 170-180.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(169).Insert()

$ws.Cells.Item(169, 1).Value2  = 4
$ws.Cells.Item(169, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(169, 3).Value2  = "Los Lagos"
$ws.Cells.Item(169, 4).Value2  = 45223
$ws.Cells.Item(169, 5).Value2  = 10
$ws.Cells.Item(169, 6).Value2  = 100112031
$ws.Cells.Item(169, 7).Value2  = "Poroto verde"
$ws.Cells.Item(169, 8).Value2  = "Magnum"
$ws.Cells.Item(169, 9).Value2  = "Primera"
$ws.Cells.Item(169, 10).Value2 = 45
$ws.Cells.Item(169, 11).Value2 = 37000
$ws.Cells.Item(169, 12).Value2 = 37000
$ws.Cells.Item(169, 13).Value2 = 37000
$ws.Cells.Item(169, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(169, 15).Value2 = "Perú"
$ws.Cells.Item(169, 16).Value2 = 1480
$ws.Cells.Item(169, 17).Value2 = 25
$ws.Cells.Item(169, 18).Value2 = "Hortaliza"
